# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: advance the list date in A1 by one day (17 Jan 2024 -> 18 Jan 2024)
$ws.Range("A1").Value = 45309

# Step 2: update the price column (D30:D33) with the new prices
$ws.Range("D30").Value = 1956.522
$ws.Range("D31").Value = 2316.776
$ws.Range("D32").Value = 2616.068
$ws.Range("D33").Value = 3297.817
